$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new column B (shifts current LOCATION column B -> C) ---
$ws.Columns.Item(2).Insert()

# --- Step 2: insert new rows for the vegetables_df2-only items ---
# Current layout after column insert (rows still 1-6):
#   Row1: FIELD, (empty), LOCATION
#   Row2: carrot, (empty), In vegetables_df1 only
#   Row3: cucumber, (empty), In vegetables_df1 only
#   Row4: potato, (empty), In vegetables_df1 only
#   Row5: spinach, (empty), In vegetables_df1 only
#   Row6: tomato, (empty), In both DataFrames
#
# Target layout needs 12 rows total, so insert a row above row2 (new broccoli row)
# and 5 rows above what will become the potato row, to host green pepper x2, kale, onion, peas.

$ws.Range("A2:C2").EntireRow.Insert()
$ws.Range("A2:C2").ClearFormats()

$ws.Range("A5:C9").EntireRow.Insert()
$ws.Range("A5:C9").ClearFormats()

# --- Step 3: write header row ---
$ws.Range("A1").Value = "vegetables_df1"
$ws.Range("B1").Value = "vegetables_df2"
$ws.Range("C1").Value = "LOCATION"

# --- Step 4: write data rows to match the final target table ---
$data = @(
    @("", "broccoli", "In vegetables_df2 only"),
    @("carrot", "", "In vegetables_df1 only"),
    @("cucumber", "", "In vegetables_df1 only"),
    @("", "green pepper", "In vegetables_df2 only"),
    @("", "green pepper", "In vegetables_df2 only"),
    @("", "kale", "In vegetables_df2 only"),
    @("", "onion", "In vegetables_df2 only"),
    @("", "peas", "In vegetables_df2 only"),
    @("potato", "", "In vegetables_df1 only"),
    @("spinach", "", "In vegetables_df1 only"),
    @("tomato", "tomato", "In both DataFrames")
)

$r = 2
foreach ($row in $data) {
    if ($row[0] -ne "") {
        $ws.Cells.Item($r, 1).Value = $row[0]
    }
    if ($row[1] -ne "") {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
